$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Pgf"
$ws.Range("C2").Value = "Flt1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 2.384646333333333
$ws.Range("H2").Value = 7.153938999999999
$ws.Range("I2").Value = 0.3245205637288701
$ws.Range("J2").Value = 0.3245205637288701
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 149.829178
$ws.Range("N2").Value = 449.487534
$ws.Range("O2").Value = 0.965236887286734
$ws.Range("P2").Value = 0.965236887286734
$ws.Range("Q2").Value = 357.2895999440473
$ws.Range("R2").Value = 3215.606399496426
$ws.Range("S2").Value = 0.3132392187941907
$ws.Range("T2").Value = 0.3132392187941908

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pgf"
$ws.Range("C3").Value = "Flt1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 2.384646333333333
$ws.Range("H3").Value = 7.153938999999999
$ws.Range("I3").Value = 0.3245205637288701
$ws.Range("J3").Value = 0.3245205637288701
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.021452666666667
$ws.Range("N3").Value = 3.064358
$ws.Range("O3").Value = 0.006580452523633729
$ws.Range("P3").Value = 0.006580452523633729
$ws.Range("Q3").Value = 2.435803356240222
$ws.Range("R3").Value = 21.922230206162
$ws.Range("S3").Value = 0.002135492162560683
$ws.Range("T3").Value = 0.002135492162560684

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Pgf"
$ws.Range("C4").Value = "Flt1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 2.384646333333333
$ws.Range("H4").Value = 7.153938999999999
$ws.Range("I4").Value = 0.3245205637288701
$ws.Range("J4").Value = 0.3245205637288701
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.374661666666667
$ws.Range("N4").Value = 13.123985
$ws.Range("O4").Value = 0.02818266018963228
$ws.Range("P4").Value = 0.02818266018963228
$ws.Range("Q4").Value = 10.43202090299055
$ws.Range("R4").Value = 93.888188126915
$ws.Range("S4").Value = 0.009145852772118653
$ws.Range("T4").Value = 0.009145852772118653

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Pgf"
$ws.Range("C5").Value = "Flt1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.736173
$ws.Range("H5").Value = 11.208519
$ws.Range("I5").Value = 0.5084464522895362
$ws.Range("J5").Value = 0.5084464522895361
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 149.829178
$ws.Range("N5").Value = 449.487534
$ws.Range("O5").Value = 0.965236887286734
$ws.Range("P5").Value = 0.965236887286734
$ws.Range("Q5").Value = 559.7877294557941
$ws.Range("R5").Value = 5038.089565102146
$ws.Range("S5").Value = 0.4907712709599348
$ws.Range("T5").Value = 0.4907712709599347

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Pgf"
$ws.Range("C6").Value = "Flt1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.736173
$ws.Range("H6").Value = 11.208519
$ws.Range("I6").Value = 0.5084464522895362
$ws.Range("J6").Value = 0.5084464522895361
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.021452666666667
$ws.Range("N6").Value = 3.064358
$ws.Range("O6").Value = 0.006580452523633729
$ws.Range("P6").Value = 0.006580452523633729
$ws.Range("Q6").Value = 3.816323873978001
$ws.Range("R6").Value = 34.34691486580201
$ws.Range("S6").Value = 0.003345807740101295
$ws.Range("T6").Value = 0.003345807740101294

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Pgf"
$ws.Range("C7").Value = "Flt1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.736173
$ws.Range("H7").Value = 11.208519
$ws.Range("I7").Value = 0.5084464522895362
$ws.Range("J7").Value = 0.5084464522895361
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.374661666666667
$ws.Range("N7").Value = 13.123985
$ws.Range("O7").Value = 0.02818266018963228
$ws.Range("P7").Value = 0.02818266018963228
$ws.Range("Q7").Value = 16.344492803135
$ws.Range("R7").Value = 147.100435228215
$ws.Range("S7").Value = 0.01432937358950008
$ws.Range("T7").Value = 0.01432937358950008

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Pgf"
$ws.Range("C8").Value = "Flt1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.227394
$ws.Range("H8").Value = 3.682182
$ws.Range("I8").Value = 0.1670329839815937
$ws.Range("J8").Value = 0.1670329839815937
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 149.829178
$ws.Range("N8").Value = 449.487534
$ws.Range("O8").Value = 0.965236887286734
$ws.Range("P8").Value = 0.965236887286734
$ws.Range("Q8").Value = 183.899434102132
$ws.Range("R8").Value = 1655.094906919188
$ws.Range("S8").Value = 0.1612263975326084
$ws.Range("T8").Value = 0.1612263975326084

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Pgf"
$ws.Range("C9").Value = "Flt1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.227394
$ws.Range("H9").Value = 3.682182
$ws.Range("I9").Value = 0.1670329839815937
$ws.Range("J9").Value = 0.1670329839815937
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.021452666666667
$ws.Range("N9").Value = 3.064358
$ws.Range("O9").Value = 0.006580452523633729
$ws.Range("P9").Value = 0.006580452523633729
$ws.Range("Q9").Value = 1.253724874350667
$ws.Range("R9").Value = 11.283523869156
$ws.Range("S9").Value = 0.001099152620971751
$ws.Range("T9").Value = 0.001099152620971751

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Pgf"
$ws.Range("C10").Value = "Flt1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.227394
$ws.Range("H10").Value = 3.682182
$ws.Range("I10").Value = 0.1670329839815937
$ws.Range("J10").Value = 0.1670329839815937
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.374661666666667
$ws.Range("N10").Value = 13.123985
$ws.Range("O10").Value = 0.02818266018963228
$ws.Range("P10").Value = 0.02818266018963228
$ws.Range("Q10").Value = 5.369433481696668
$ws.Range("R10").Value = 48.32490133527001
$ws.Range("S10").Value = 0.004707433828013548
$ws.Range("T10").Value = 0.004707433828013548
